# Eindpaper_Kerncentrales_JortSiemes_V2.docx edit
#
# 1. Replace the "Vergelijking stroomproductie:" bullet list (translated
#    English/mixed energy-source names) with the Dutch percentage list,
#    dropping the now-empty spare paragraph.
# 2. Move the <w:lastRenderedPageBreak/> marker from the "A. Samenvatting
#    van belangrijkste resultaten" run to the "B. Aanbevelingen voor
#    nieuwsmedia en beleidsmakers" run.

$d = $word.ActiveDocument

function Find-ParagraphIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.TrimEnd("`r`a") -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# Part 1: "Vergelijking stroomproductie:" list -> Dutch percentage list
# ---------------------------------------------------------------------

$firstIdx = Find-ParagraphIndexByText $d "Kerncentrale"
$lastIdx  = Find-ParagraphIndexByText $d "solar"

$pFirst = $d.Paragraphs.Item($firstIdx)
# The two blank paragraphs right after "solar" belong to this block too
# (they disappear - replaced by a single trailing blank bold paragraph).
$pLast  = $d.Paragraphs.Item($lastIdx + 2)

$listRange = $d.Range($pFirst.Range.Start, $pLast.Range.End)

$newListXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>Kerncentrale</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve"> = </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>92.5%</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:t>Aardwarmte = 74,3%</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:t>Aardgas = 56,6%</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:t>Waterkracht = 41,5%</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:t>Steenkool = 40,2%</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:t>Wind = 35,4%</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:b/>
      <w:bCs/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>Zonne-</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>energie</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve"> = 24,9%</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:b/>
      <w:bCs/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@

$listRange.InsertXML($newListXml)

# ---------------------------------------------------------------------
# Part 2: move <w:lastRenderedPageBreak/> from "A. Samenvatting..." run
#         to "B. Aanbevelingen..." run.
# ---------------------------------------------------------------------

$idxA = Find-ParagraphIndexByText $d "A. Samenvatting van belangrijkste resultaten"
$idxB = $idxA + 1

$pA = $d.Paragraphs.Item($idxA)
$pB = $d.Paragraphs.Item($idxB)

$pageBreakRange = $d.Range($pA.Range.Start, $pB.Range.End)

$pageBreakXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:b/>
      <w:bCs/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:b/>
      <w:bCs/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:t>A. Samenvatting van belangrijkste resultaten</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:b/>
      <w:bCs/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Civil Premium" w:hAnsi="Civil Premium"/>
      <w:b/>
      <w:bCs/>
      <w:lang w:val="nl-NL"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>B. Aanbevelingen voor nieuwsmedia en beleidsmakers</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
'@

$pageBreakRange.InsertXML($pageBreakXml)
